# Finished bookstore part 1
# Update daily-change percentages (I3:I6) and price/weight figures
# (G3:G7, I7) on the "AwesomeSheet" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PercentText($cell, $text) {
    # Cells I3:I6 store their "Daily Change %" as literal text (shared
    # string) even though the cell is number-formatted as a percentage.
    # Assigning a "-x.xx%" string straight to .Value would get parsed as
    # a numeric percentage instead of preserved as text, so we briefly
    # switch the cell to a text format, write the literal string, then
    # restore the original percentage format.
    $range = $ws.Range($cell)
    $originalFormat = $range.NumberFormat
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.NumberFormat = $originalFormat
}

# Daily Change % column (stored as shared-string text)
Set-PercentText "I3" "-0.60%"
Set-PercentText "I4" "-1.16%"
Set-PercentText "I5" "-3.24%"
Set-PercentText "I6" "-1.92%"

# Prev Close column (numeric values)
$ws.Range("G3").Value = 145.16
$ws.Range("G4").Value = 136.67
$ws.Range("G5").Value = 8.34
$ws.Range("G6").Value = 23.99
$ws.Range("G7").Value = 369.885

# Daily Change % for row 7 (stored as a raw number, not text)
$ws.Range("I7").Value = 0.1254
